# edit.ps1 — apply "this_month_production.xlsx" changes via Excel COM-interop
#
# Summary of the edit:
#  - Remove the "المجموعة" (Group) and "الوردية" (Shift) columns (old C:D).
#  - Remove the second employee row (old row 3).
#  - Replace the remaining employee's data with a new employee record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old row 3 (second employee) entirely.
$ws.Rows.Item(3).Delete()

# 2) Drop the "المجموعة" (old C) and "الوردية" (old D) columns, shifting
#    everything after them (الانتاج الكلي, اليوم 1..31) one and two slots left.
$ws.Range("C:D").Delete()

# 3) Overwrite the remaining employee row with the new record.
$ws.Range("A2").Value = 185
$ws.Range("B2").Value = "عبد المنعم محمد القضيب"
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 4000
$ws.Range("E2").Value = 4000

# 4) Clear any leftover day values beyond day 2 (old "اليوم 27" data used to
#    live here before the column shift).
$ws.Range("F2:AH2").ClearContents()
